# "refactor with sew to implement code chunks"
#
# knitxl now renders the R source of each code chunk above its output, in
# a monospace font. For this snapshot, the chunk
#   setNames(LETTERS[1:3], c('a', 'b', 'c'))
# gets written into A1 (Courier New), a blank spacer row is left at row 2,
# and the pre-existing a/A/b/B/c/C value grid (previously rows 1-3) shifts
# down to occupy rows 3-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the existing grid (rows 1-3, columns A-B) before touching
# anything, since row 1 is about to be overwritten with the chunk source.
$vals = @()
for ($r = 1; $r -le 3; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 2; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $vals += ,$rowVals
}

# Clear the original 3x2 block; it will be rebuilt two rows lower.
$ws.Range("A1:B3").Clear()

# Row 1: the chunk's source code, rendered in a monospace font.
$ws.Cells.Item(1, 1).Value = "setNames(LETTERS[1:3], c('a', 'b', 'c'))"
$ws.Cells.Item(1, 1).Font.Name = "Courier New"

# Row 2 is left empty (spacer between code and its output).

# Rows 3-5: the original value grid, shifted down by two rows.
for ($i = 0; $i -lt 3; $i++) {
    $targetRow = 3 + $i
    $ws.Cells.Item($targetRow, 1).Value = $vals[$i][0]
    $ws.Cells.Item($targetRow, 2).Value = $vals[$i][1]
}

$wb.Save()
